{"js": "// Rename the \"intermediate\" tables/commands so CREATE TABLE / COPY statements\n// target the \"<name>_ intermedio\" table instead of \"<name>\" directly.\n// (First \"actores_pelis\" CREATE TABLE keeps no spaces: \"actores_pelis_intermedio\".)\nconst replacements = [\n  [\"CREATE TABLE actores_pelis(\", \"CREATE TABLE actores_pelis_intermedio(\"],\n  [\"CREATE TABLE actores(\", \"CREATE TABLE actores_ intermedio (\"],\n  [\"CREATE TABLE directores(\", \"CREATE TABLE directores_ intermedio (\"],\n  [\"CREATE TABLE peliculas(\", \"CREATE TABLE pel\u00edculas_ intermedio (\"],\n  [\"CREATE TABLE pelis_comentarios(\", \"CREATE TABLE pelis_comentarios_ intermedio (\"],\n  [\"CREATE TABLE pelis_directores(\", \"CREATE TABLE pelis_directores_ intermedio (\"],\n  [\"COPY actores_pelis(\", \"COPY actores_pelis_ intermedio (\"],\n  [\"COPY actores(\", \"COPY actores_ intermedio (\"],\n  [\"COPY directores(\", \"COPY directores_ intermedio (\"],\n  [\"COPY peliculas(\", \"COPY pel\u00edculas_ intermedio (\"],\n  [\"COPY pelis_comentarios (\", \"COPY pelis_comentarios_ intermedio (\"],\n  [\"COPY pelis_directores (\", \"COPY pelis_directores_ intermedio (\"],\n];\n\nconst body = context.document.body;\n\nfor (const [needle, replacement] of replacements) {\n  const results = body.search(needle, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  if (results.items.length === 0) {\n    continue;\n  }\n\n  // Each needle is unique in this document, so take the first (only) hit.\n  results.items[0].insertText(replacement, \"Replace\");\n  await context.sync();\n}\n", "ps1": "# Rename the \"intermediate\" tables/commands so CREATE TABLE / COPY statements\n# target the \"<name>_ intermedio\" table instead of \"<name>\" directly.\n# (First \"actores_pelis\" CREATE TABLE keeps no spaces: \"actores_pelis_intermedio\".)\n\n$d = $word.ActiveDocument\n\n$replacements = @(\n    @(\"CREATE TABLE actores_pelis(\", \"CREATE TABLE actores_pelis_intermedio(\"),\n    @(\"CREATE TABLE actores(\", \"CREATE TABLE actores_ intermedio (\"),\n    @(\"CREATE TABLE directores(\", \"CREATE TABLE directores_ intermedio (\"),\n    @(\"CREATE TABLE peliculas(\", \"CREATE TABLE pel\u00edculas_ intermedio (\"),\n    @(\"CREATE TABLE pelis_comentarios(\", \"CREATE TABLE pelis_comentarios_ intermedio (\"),\n    @(\"CREATE TABLE pelis_directores(\", \"CREATE TABLE pelis_directores_ intermedio (\"),\n    @(\"COPY actores_pelis(\", \"COPY actores_pelis_ intermedio (\"),\n    @(\"COPY actores(\", \"COPY actores_ intermedio (\"),\n    @(\"COPY directores(\", \"COPY directores_ intermedio (\"),\n    @(\"COPY peliculas(\", \"COPY pel\u00edculas_ intermedio (\"),\n    @(\"COPY pelis_comentarios (\", \"COPY pelis_comentarios_ intermedio (\"),\n    @(\"COPY pelis_directores (\", \"COPY pelis_directores_ intermedio (\")\n)\n\nforeach ($pair in $replacements) {\n    $needle = $pair[0]\n    $replacement = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Execute($needle, $false, $false, $false, $false, $false, $true, 1, $false, $replacement, 2)\n}\n\nWrite-Output \"done\"\n"}
